$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextValue 'D2' '59.132.05'
Set-TextValue 'E2' '  -2.33%  '
Set-TextValue 'D3' '2.570.39'
Set-TextValue 'E3' '  -2.61%  '
Set-TextValue 'E4' '  -0.02%  '
Set-TextValue 'D5' '554.24'
Set-TextValue 'E5' '  -2.29%  '
Set-TextValue 'D6' '141.30'
Set-TextValue 'E6' '  -3.69%  '
Set-TextValue 'E7' '  +0.37%  '
Set-TextValue 'E8' '  -2.18%  '
Set-TextValue 'D9' '2.574.17'
Set-TextValue 'E9' '  -3.37%  '
Set-TextValue 'D10' '6.65'
Set-TextValue 'E10' '  -2.90%  '
Set-TextValue 'E11' '  -1.54%  '
Set-TextValue 'E12' '  +11.87%  '
Set-TextValue 'D13' '0.351'
Set-TextValue 'E13' '  +2.38%  '
Set-TextValue 'D14' '3.023.14'
Set-TextValue 'E14' '  -2.62%  '
Set-TextValue 'D15' '59.167.94'
Set-TextValue 'E15' '  -2.29%  '
Set-TextValue 'D16' '23.00'
Set-TextValue 'E16' '  +4.47%  '
Set-TextValue 'E17' '  -1.50%  '
Set-TextValue 'D18' '2.576.34'
Set-TextValue 'E18' '  -2.89%  '
Set-TextValue 'E19' '  -0.10%  '
Set-TextValue 'D20' '335.58'
Set-TextValue 'E20' '  -2.29%  '
Set-TextValue 'D21' '10.29'
Set-TextValue 'E21' '  -1.48%  '
Set-TextValue 'E22' '  +0.01%  '
Set-TextValue 'E23' '  -0.24%  '
Set-TextValue 'B24' 'Polygon'
Set-TextValue 'C24' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 'D24' '0.473'
Set-TextValue 'E24' '  +7.24%  '
Set-TextValue 'B25' 'Litecoin'
Set-TextValue 'C25' 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue 'D25' '62.73'
Set-TextValue 'E25' '  -4.95%  '
Set-TextValue 'D26' '0.999'
Set-TextValue 'E26' '  +0.61%  '
Set-TextValue 'D27' '0.158'
Set-TextValue 'E27' '  -3.98%  '
Set-TextValue 'D28' '7.37'
Set-TextValue 'E28' '  -0.56%  '
Set-TextValue 'D29' '0.0₃0770'
Set-TextValue 'E29' '  -4.76%  '
Set-TextValue 'E30' '  +0.11%  '
Set-TextValue 'D31' '6.15'
Set-TextValue 'E31' '  -1.01%  '
Set-TextValue 'E32' '  -3.08%  '
Set-TextValue 'D33' '157.55'
Set-TextValue 'E33' '  -1.02%  '
Set-TextValue 'D34' '18.97'
Set-TextValue 'E34' '  -1.38%  '
Set-TextValue 'E35' '  -2.15%  '
Set-TextValue 'B36' 'ImmutableX'
Set-TextValue 'C36' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D36' '1.16'
Set-TextValue 'E36' '  -0.63%  '
Set-TextValue 'B37' 'Fetch.AI'
Set-TextValue 'C37' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 'D37' '0.892'
Set-TextValue 'E37' '  -0.99%  '
Set-TextValue 'E38' '  -0.58%  '
Set-TextValue 'E39' '  -6.08%  '
Set-TextValue 'E40' '  -2.78%  '
Set-TextValue 'E41' '  -0.18%  '
Set-TextValue 'D42' '289.84'
Set-TextValue 'E42' '  -3.92%  '
Set-TextValue 'D43' '135.74'
Set-TextValue 'E43' '  +5.06%  '
Set-TextValue 'E44' '  +0.48%  '
Set-TextValue 'E45' '  -1.55%  '
Set-TextValue 'D46' '0.589'
Set-TextValue 'E46' '  -2.61%  '
Set-TextValue 'E47' '  -0.36%  '
Set-TextValue 'D48' '0.0529'
Set-TextValue 'E48' '  -3.05%  '
Set-TextValue 'D49' '0.0232'
Set-TextValue 'E49' '  -2.11%  '
Set-TextValue 'D50' '18.57'
Set-TextValue 'E50' '  -0.65%  '
Set-TextValue 'D51' '1.941.65'
Set-TextValue 'E51' '  -0.98%  '
